$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 516.6667
$ws.Range("J9").Value = 150
$ws.Range("L9").Value = 150
$ws.Range("N9").Value = -488
$ws.Range("H32").Value = 9406.25
$ws.Range("I32").Value = 7001
$ws.Range("J32").Value = 9749.857
$ws.Range("K32").Value = 7001
$ws.Range("L32").Value = 9749.857
$ws.Range("M32").Value = -6675
$ws.Range("N32").Value = -10401.857
$ws.Range("H40").Value = 6334.7827
$ws.Range("I40").Value = 1100
$ws.Range("J40").Value = 6572.727
$ws.Range("K40").Value = 1100
$ws.Range("L40").Value = 6572.727
$ws.Range("M40").Value = -925
$ws.Range("N40").Value = -6922.727
$ws.Range("H48").Value = 465.75
$ws.Range("I48").Value = 387.66666
$ws.Range("J48").Value = 700
$ws.Range("K48").Value = 1162.99998
$ws.Range("L48").Value = 2100
$ws.Range("M48").Value = -870.9999800000001
$ws.Range("N48").Value = -2684
$ws.Range("H56").Value = 465.75
$ws.Range("I56").Value = 387.66666
$ws.Range("J56").Value = 700
$ws.Range("K56").Value = 1162.99998
$ws.Range("L56").Value = 2100
$ws.Range("M56").Value = -628.9999800000001
$ws.Range("N56").Value = -3168
$ws.Range("H62").Value = 2666
$ws.Range("I62").Value = 2999
$ws.Range("K62").Value = 2999
$ws.Range("M62").Value = -2375
$ws.Range("H65").Value = 2666
$ws.Range("I65").Value = 2999
$ws.Range("K65").Value = 14995
$ws.Range("M65").Value = -11875
$ws.Range("H74").Value = 5400
$ws.Range("I74").Value = 4666.6665
$ws.Range("K74").Value = 4666.6665
$ws.Range("M74").Value = -3730.6665
$ws.Range("H77").Value = 5400
$ws.Range("I77").Value = 4666.6665
$ws.Range("K77").Value = 23333.3325
$ws.Range("M77").Value = -18653.3325
$ws.Range("H86").Value = 1799.1428
$ws.Range("I86").Value = 869
$ws.Range("J86").Value = 4124.5
$ws.Range("K86").Value = 869
$ws.Range("L86").Value = 4124.5
$ws.Range("M86").Value = 254
$ws.Range("N86").Value = -6370.5
$ws.Range("H89").Value = 1799.1428
$ws.Range("I89").Value = 869
$ws.Range("J89").Value = 4124.5
$ws.Range("K89").Value = 4345
$ws.Range("L89").Value = 20622.5
$ws.Range("M89").Value = 1271
$ws.Range("N89").Value = -31854.5
$ws.Range("H107").Value = 2307
$ws.Range("I107").Value = 1530
$ws.Range("K107").Value = 1530
$ws.Range("M107").Value = 390
$ws.Range("H113").Value = 3388.6667
$ws.Range("I113").Value = 3474.5
$ws.Range("K113").Value = 3474.5
$ws.Range("M113").Value = -220.5
$ws.Range("H132").Value = 7261.5
$ws.Range("I132").Value = 6530.3335
$ws.Range("J132").Value = 8577.6
$ws.Range("K132").Value = 19591.0005
$ws.Range("L132").Value = 25732.8
$ws.Range("M132").Value = -17061.0005
$ws.Range("N132").Value = -30792.8
$ws.Range("H137").Value = 1689.5
$ws.Range("I137").Value = 1545.75
$ws.Range("J137").Value = 1977
$ws.Range("K137").Value = 4637.25
$ws.Range("L137").Value = 5931
$ws.Range("M137").Value = -2087.25
$ws.Range("N137").Value = -11031

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5205.875
$ws.Range("I32").Value = 5205.875
$ws.Range("K32").Value = 5205.875
$ws.Range("M32").Value = -4918.875
$ws.Range("H36").Value = 3006.5
$ws.Range("I36").Value = 3006.5
$ws.Range("K36").Value = 3006.5
$ws.Range("M36").Value = -2660.5
$ws.Range("H61").Value = 3975.8333
$ws.Range("I61").Value = 3975.8333
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3975.8333
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -3763.8333
$ws.Range("H110").Value = 877.75
$ws.Range("I110").Value = 670.3333
$ws.Range("K110").Value = 670.3333
$ws.Range("M110").Value = 1374.6667
$ws.Range("H122").Value = 9995.817999999999
$ws.Range("I122").Value = 9995.817999999999
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 29987.454
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -27537.454
$ws.Range("H132").Value = 1266.6666
$ws.Range("J132").Value = 1266.6666
$ws.Range("L132").Value = 3799.9998
$ws.Range("N132").Value = -8859.9998
$ws.Range("H136").Value = 3975.8333
$ws.Range("I136").Value = 3975.8333
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 11927.4999
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -9377.499899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 587.5
$ws.Range("J16").Value = 600
$ws.Range("L16").Value = 600
$ws.Range("N16").Value = -1174
$ws.Range("H99").Value = 2999.5
$ws.Range("I99").Value = 2999
$ws.Range("K99").Value = 2999
$ws.Range("M99").Value = -1501
$ws.Range("H105").Value = 1375.6
$ws.Range("J105").Value = 1000
$ws.Range("L105").Value = 1000
$ws.Range("N105").Value = -4494
$ws.Range("H107").Value = 800
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H113").Value = 587.5
$ws.Range("J113").Value = 600
$ws.Range("L113").Value = 600
$ws.Range("N113").Value = -4940
$ws.Range("H126").Value = 2999.5
$ws.Range("I126").Value = 2999
$ws.Range("K126").Value = 8997
$ws.Range("M126").Value = -6527

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1250492.2
$ws.Range("I4").Value = 1250492.2
$ws.Range("K4").Value = 3751476.6
$ws.Range("M4").Value = -3751364.6
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 555
$ws.Range("I19").Value = 555
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 555
$ws.Range("L19").ClearContents()
$ws.Range("M19").Value = -267
$ws.Range("N19").Value = 0
$ws.Range("H102").Value = 2999.5
$ws.Range("I102").Value = 2999.5
$ws.Range("K102").Value = 2999.5
$ws.Range("M102").Value = -1377.5
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = 0
$ws.Range("H132").Value = 3307.5
$ws.Range("I132").Value = 3307.5
$ws.Range("K132").Value = 9922.5
$ws.Range("M132").Value = -7392.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 25000000
$ws.Range("I3").Value = 25000000
$ws.Range("K3").Value = 25000000
$ws.Range("M3").Value = -24999888
$ws.Range("H7").Value = 5500
$ws.Range("I7").Value = 5500
$ws.Range("K7").Value = 5500
$ws.Range("M7").Value = -5388
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H15").Value = 25000000
$ws.Range("I15").Value = 25000000
$ws.Range("K15").Value = 25000000
$ws.Range("M15").Value = -24999830
$ws.Range("H40").Value = 7883.5
$ws.Range("I40").Value = 7883.5
$ws.Range("K40").Value = 7883.5
$ws.Range("M40").Value = -7747.5
$ws.Range("H43").Value = 9402.799999999999
$ws.Range("I43").Value = 9000
$ws.Range("J43").Value = 10007
$ws.Range("K43").Value = 9000
$ws.Range("L43").Value = 10007
$ws.Range("M43").Value = -8807
$ws.Range("N43").Value = -10393
$ws.Range("H100").Value = 3579.6
$ws.Range("I100").Value = 3849.5
$ws.Range("K100").Value = 3849.5
$ws.Range("M100").Value = -3308.5
$ws.Range("H122").Value = 4999.25
$ws.Range("I122").Value = 4999.25
$ws.Range("K122").Value = 14997.75
$ws.Range("M122").Value = -12547.75
$ws.Range("H126").Value = 5500
$ws.Range("I126").Value = 5500
$ws.Range("K126").Value = 16500
$ws.Range("M126").Value = -14030
$ws.Range("H132").Value = 4538.5835
$ws.Range("I132").Value = 4396.3
$ws.Range("K132").Value = 13188.9
$ws.Range("M132").Value = -10658.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 3000000
$ws.Range("I3").Value = 3000000
$ws.Range("K3").Value = 3000000
$ws.Range("M3").Value = -2999886
$ws.Range("H32").Value = 3000
$ws.Range("I32").Value = 3000
$ws.Range("K32").Value = 3000
$ws.Range("M32").Value = -2683
$ws.Range("H122").Value = 7599.4
$ws.Range("I122").Value = 6332.3335
$ws.Range("K122").Value = 18997.0005
$ws.Range("M122").Value = -16547.0005
